$d = $word.ActiveDocument

# The document's last paragraph is the trailing empty paragraph that
# carries the "_GoBack" bookmark (right before the final <w:sectPr>).
# Append two new paragraphs right after it:
#   1. An "Update" heading (Heading 2 / style "21").
#   2. A body paragraph describing the new EquationAssigner class,
#      using manual line breaks (Word's vertical-tab char, Chr(11))
#      between the individual lines so they land in a single run,
#      just like text typed with Shift+Enter in Word.

$lastIndex = $d.Paragraphs.Count
$lastRange = $d.Paragraphs($lastIndex).Range
$lastRange.InsertParagraphAfter()

# --- Paragraph 1: "Update" heading ---
$headingPara = $d.Paragraphs($lastIndex + 1)
$headingPara.Style = "21"
$headingPara.Range.Text = "🔄 Update: New Class – EquationAssigner"

# Make room for the second new paragraph.
$headingPara.Range.InsertParagraphAfter()

# --- Paragraph 2: EquationAssigner class description ---
$bodyPara = $d.Paragraphs($lastIndex + 2)
$bodyPara.Style = "a1"

$lineBreak = [char]11
$bodyLines = @(
    "class EquationAssigner:",
    "  + assign_equation(df: pd.DataFrame) -> pd.DataFrame",
    "  + compute_manning(df: pd.DataFrame) -> pd.DataFrame",
    "",
    "🔗 Assigns governing equations to each spatial cell and estimates Manning's coefficient where needed.",
    "To be integrated before training, after feature engineering."
)
$bodyPara.Range.Text = [string]::Join($lineBreak, $bodyLines)
